$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "MaxDFR" in L1
$ws.Range("L1").Value = "MaxDFR"

# Add new value "5000" in L2 (stored as text, matching quotePrefix style of the row's data cells)
$ws.Range("L2").Style = $ws.Range("J2").Style
$ws.Range("L2").Value = "'5000"

# Column K's content now drives a best-fit column width (Excel
# recalculates best-fit widths for the sheet as data changes)
$ws.Columns.Item(11).ColumnWidth = 22.7

# Update selection to reflect the newly active cell after edit
$ws.Range("L1").Select()
